# Update countries & provincias Spain
# Applies updated statistics for a handful of country rows on the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Map of row number -> new values for columns B..H
# (A column, the country id/name, stays unchanged)
$updates = @{
    4   = @(970757, 10106, 118633, 797183, 15116, 685, 54941)
    8   = @(157114, 601,   109800, 41430,  2570,  7,   5884)
    10  = @(110130, 2357,  29140,  78185,  1776,  99,  2805)
    32  = @(13201,  478,   2936,   9993,   111,   3,   272)
    55  = @(4065,   168,   593,    3311,   1,     2,   161)
    106 = @(505,    53,    120,    378,    2,     0,   7)
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Column B is index 2
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
